# Updates the "cryptos" price list (Price / Volume(1h) columns, plus one
# ranking swap between ARBITRUM and Kaspa in rows 38-39) to the latest
# scrape snapshot, mirroring the GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking strings
# (e.g. "301.40", "19.00", "0.100") are preserved exactly as text, matching
# the original inline-string cells, instead of being coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "42.771.31"
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("D3").Value = "2.309.12"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "301.40"
$ws.Range("E5").Value = "  -0.34%  "

$ws.Range("D6").Value = "95.22"
$ws.Range("E6").Value = "  -0.86%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "0.491"
$ws.Range("E9").Value = "  -0.93%  "

$ws.Range("D10").Value = "34.13"
$ws.Range("E10").Value = "  -1.86%  "

$ws.Range("D11").Value = "19.00"
$ws.Range("E11").Value = "  +2.13%  "

$ws.Range("D12").Value = "0.0782"
$ws.Range("E12").Value = "  +0.12%  "

$ws.Range("E13").Value = "  -0.01%  "

$ws.Range("E14").Value = "  -1.42%  "

$ws.Range("D15").Value = "2.673.54"
$ws.Range("E15").Value = "  +0.82%  "

$ws.Range("D16").Value = "2.323.00"
$ws.Range("E16").Value = "  +1.04%  "

$ws.Range("E17").Value = "  +1.72%  "

$ws.Range("D18").Value = "42.717.24"
$ws.Range("E18").Value = "  +0.28%  "

$ws.Range("D19").Value = "12.18"
$ws.Range("E19").Value = "  -4.91%  "

$ws.Range("D20").Value = "6.12"
$ws.Range("E20").Value = "  +2.24%  "

$ws.Range("D21").Value = "0.0₃0889"
$ws.Range("E21").Value = "  -0.25%  "

$ws.Range("D22").Value = "67.71"
$ws.Range("E22").Value = "  +1.00%  "

$ws.Range("D23").Value = "2.27"
$ws.Range("E23").Value = "  +7.34%  "

$ws.Range("D24").Value = "234.94"
$ws.Range("E24").Value = "  -0.38%  "

$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("E26").Value = "  +1.05%  "

$ws.Range("D27").Value = "24.25"
$ws.Range("E27").Value = "  -1.16%  "

$ws.Range("E28").Value = "  +15.09%  "

$ws.Range("D29").Value = "166.40"
$ws.Range("E29").Value = "  -0.90%  "

$ws.Range("D30").Value = "9.13"
$ws.Range("E30").Value = "  +1.71%  "

$ws.Range("D31").Value = "32.05"
$ws.Range("E31").Value = "  -2.27%  "

$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("E33").Value = "  +0.95%  "

$ws.Range("D34").Value = "17.67"
$ws.Range("E34").Value = "  -0.17%  "

$ws.Range("D35").Value = "4.43"
$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("D36").Value = "0.0695"
$ws.Range("E36").Value = "  +1.76%  "

$ws.Range("D37").Value = "2.33"
$ws.Range("E37").Value = "  -0.85%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.100"
$ws.Range("E38").Value = "  +0.04%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "1.77"
$ws.Range("E39").Value = "  +1.85%  "

$ws.Range("D40").Value = "2.71"
$ws.Range("E40").Value = "  +1.18%  "

$ws.Range("D41").Value = "0.108"
$ws.Range("E41").Value = "  -0.57%  "

$ws.Range("D42").Value = "20.65"
$ws.Range("E42").Value = "  +13.33%  "

$ws.Range("D43").Value = "1.924.79"

$ws.Range("D44").Value = "0.0278"
$ws.Range("E44").Value = "  -0.31%  "

$ws.Range("D45").Value = "10.15"
$ws.Range("E45").Value = "  +0.19%  "

$ws.Range("E46").Value = "  -2.31%  "

$ws.Range("D47").Value = "2.73"
$ws.Range("E47").Value = "  -1.11%  "

$ws.Range("E48").Value = "  +2.24%  "

$ws.Range("D49").Value = "2.540.84"
$ws.Range("E49").Value = "  +0.88%  "

$ws.Range("D50").Value = "53.20"
$ws.Range("E50").Value = "  -0.51%  "

$ws.Range("D51").Value = "72.04"
$ws.Range("E51").Value = "  +1.83%  "

# Restore the default (General) style on the touched range so no stray
# cell-format/style metadata is introduced by the text-format trick above.
$ws.Range("D2:E51").Style = "Normal"
